$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells marked "Yes" with the newly-applied Arial 14 (theme color 1) font + border style
# (matches the distinctive new style bucket introduced in the edit)
$rng = $ws.Range("C5:I5")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("O5")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C6:I6")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("O6")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C7:I7")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("O7")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C8:I8")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("O8")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C9:I9")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("O9")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C10:I10")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("O10")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C11:I11")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("O11")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C12:I12")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("O12")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C13:E13")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("O13")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C14:E14")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("O14")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C15:E15")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("H15")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C16:E16")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C17:E17")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C18:E18")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C19:D19")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("H19")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C20:D20")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("H20")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C21:D21")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C22:D22")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C23:D23")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C24:D24")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C25:D25")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C26:D26")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C27:D27")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C28:D28")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C29:D29")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C30:D30")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C31:D31")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C32:D32")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C33:D33")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C34:D34")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C35:D35")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C36:D36")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C37:D37")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C38:D38")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C39:D39")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C40:D40")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C41:D41")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C42:D42")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C43:D43")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C44:D44")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C45:D45")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C46:D46")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C47:D47")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

$rng = $ws.Range("C48:D48")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14
$rng.Font.ThemeColor = 1

# Cells marked "Yes" using the pre-existing Arial 14 font/border bucket (no new font)
$rng = $ws.Range("E20")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14

$rng = $ws.Range("E21")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14

$rng = $ws.Range("E22")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14

$rng = $ws.Range("E23")
$rng.Value = "Yes"
$rng.Font.Name = "Arial"
$rng.Font.Size = 14

# Cell(s) that already contained "Yes" but whose style bucket is normalised to the plain font
$rng = $ws.Range("E19")
$rng.Font.Name = "Arial"
$rng.Font.Size = 14

# Restore the last active selection to match the saved workbook state
$ws.Range("I13").Select()